$d = $word.ActiveDocument

# Locate the target paragraph: the one containing "Форма не работает на локалке"
$targetIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $t = $d.Paragraphs.Item($i).Range.Text
    if ($t -like "*Форма не работает на локалке*") {
        $targetIndex = $i
        break
    }
}
if ($targetIndex -eq -1) {
    throw "Target paragraph not found"
}

$p = $d.Paragraphs.Item($targetIndex)
$r = $p.Range
$insertPoint = $r.Duplicate
$insertPoint.Collapse(0)

$newXml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:t xml:space="preserve">Вот так заставляем работать плагин при помощи </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>шорткода</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">. Форма не работает на </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>локалке</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p><w:p><w:r><w:t>ХОСТИНГ!</w:t></w:r></w:p><w:p><w:r><w:t>PLACEHOLDER_HYPERLINK</w:t></w:r></w:p><w:p><w:r><w:t>На домен нет денег</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">В файле менеджеров нужно подвязать </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>новуб</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>бд</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> в файле </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>вп_конфиг</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p><w:p><w:r><w:t xml:space="preserve">После того как закинул сайт на хостинг нужна </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>рекаптча</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> от </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>ботов,</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>в</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>плагин</w:t></w:r><w:r><w:t xml:space="preserve">е </w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>contact</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>form</w:t></w:r><w:r><w:t xml:space="preserve"> 7</w:t></w:r><w:r><w:t xml:space="preserve"> в импорте </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>рекаптча</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> от </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>гугл</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">, </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>каптчу</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> 3го </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>покления</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> (нужен домен </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>е№ный</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> ;%ять)</w:t></w:r><w:r><w:t>, там можно достать ключи интеграции. Все просто….</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>
'@

$insertPoint.InsertXML($newXml)

Write-Output ("Paragraphs after InsertXML=" + $d.Paragraphs.Count)

# Locate the placeholder paragraph for the hyperlink
$hlIndex = -1
for ($i = $targetIndex; $i -le $d.Paragraphs.Count; $i++) {
    $t = $d.Paragraphs.Item($i).Range.Text
    if ($t -like "*PLACEHOLDER_HYPERLINK*") {
        $hlIndex = $i
        break
    }
}
if ($hlIndex -eq -1) {
    throw "Hyperlink placeholder paragraph not found"
}

$hp = $d.Paragraphs.Item($hlIndex)
$hr = $hp.Range
$hr.MoveEnd(1, -1)
$h = $d.Hyperlinks.Add($hr, "https://2domains.ru/manager", $null, $null, "https://2domains.ru/manager")
$h.Range.Style = "Hyperlink"

$hp2 = $d.Paragraphs.Item($hlIndex)
$hr2 = $hp2.Range
$hr2.Collapse(0)
$hr2.MoveEnd(1, -1)
$hr2.InsertAfter(" (тут хостинг)")

Write-Output ("Final paragraph count=" + $d.Paragraphs.Count)
